$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (avoids Excel's automatic date/number
# auto-detection for strings like "2024-09-12" or "23"), while leaving the
# cell's style untouched (no lingering quotePrefix / number-format changes).
function Set-TextValue($range, [string]$value) {
    if ($value -eq "") {
        # A lone leading apostrophe produces an empty TEXT cell (matches the
        # workbook's existing representation of blank text cells).
        $range.Value = "'"
    } else {
        $range.NumberFormat = "@"
        $range.Value = $value
    }
    $range.ClearFormats()
}

# Remove the four trailing rows (old rows 8-11) that no longer exist.
$ws.Range("A8:A11").EntireRow.Delete()

# Row 5 updates
$ws.Range("A5").Value = 8
Set-TextValue $ws.Range("B5") "23"
Set-TextValue $ws.Range("C5") "2024-09-12"
Set-TextValue $ws.Range("D5") ""
$ws.Range("E5").Value = 23
$ws.Range("F5").Value = 23
Set-TextValue $ws.Range("G5") "M1B"
Set-TextValue $ws.Range("H5") "ASDF"
Set-TextValue $ws.Range("I5") "ASFASDF"
Set-TextValue $ws.Range("J5") "ASDFAS"
Set-TextValue $ws.Range("K5") "DFA"

# Row 6 updates
$ws.Range("A6").Value = 11
Set-TextValue $ws.Range("B6") "234"
Set-TextValue $ws.Range("F6") ""
Set-TextValue $ws.Range("G6") "M1A"
Set-TextValue $ws.Range("H6") "Explicación"
Set-TextValue $ws.Range("I6") "safd"
Set-TextValue $ws.Range("J6") "asdf"
Set-TextValue $ws.Range("K6") "asdfas"

# Row 7 updates
$ws.Range("A7").Value = 10
Set-TextValue $ws.Range("B7") "HOLA MIBIDA"
Set-TextValue $ws.Range("C7") "2024-09-13"
Set-TextValue $ws.Range("E7") "3232"
Set-TextValue $ws.Range("F7") ""
Set-TextValue $ws.Range("G7") "TM"
Set-TextValue $ws.Range("H7") "Teórico/Práctica"
Set-TextValue $ws.Range("I7") "sadf"
Set-TextValue $ws.Range("J7") "asdf"
Set-TextValue $ws.Range("K7") "asdf"
